# Apply "Ajout synthese des chiffres sur recuit" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, fix the irregular source data:
#  - block starting at label row 265 has an extra (11th) data row at row 276 that must be
#    cleared so every block has exactly 10 data rows.
$ws.Range("A276:B276").ClearContents()

#  - block starting at label row 288 is missing its 10th data row; add it at row 298.
$ws.Range("A298").Value = 0
$ws.Range("B298").Value = 0

# Data-start row (first data row) of every one of the 27 blocks (label row + 1).
$dataStarts = @(2,13,24,35,46,57,68,79,90,101,112,123,134,145,156,167,178,189,200,211,222,233,244,255,266,278,289)

foreach ($r in $dataStarts) {
    $rEnd = $r + 9
    $ws.Range("D$r").Formula = "=AVERAGE(A$r`:A$rEnd)"
    $ws.Range("E$r").Formula = "=AVERAGE(B$r`:B$rEnd)"
}

# Highlight the first synthesis row (D2:E2) in yellow.
$ws.Range("D2:E2").Interior.Color = 65535

# Restore the view/selection state recorded in the saved workbook.
$excel.ActiveWindow.ScrollRow = 26
$ws.Range("D35").Select()

$excel.Calculate()
$wb.Save()
